$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the total "Valor Mora" amount
$ws.Range("E11").Value = 17082

# Update worker / period counts
$ws.Range("C13").Value = 4
$ws.Range("F13").Value = 3

# Insert 3 new rows above the existing data row (row 16), pushing it down to row 19.
# This keeps the original last row's formatting (incl. the bottom border) intact on row 19,
# and the newly inserted rows 16-18 inherit row 16's (now row 19's) current formatting.
$ws.Rows("16:18").Insert()
$ws.Range("B19:J19").Copy()
$ws.Range("B16:J18").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the new worker rows
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "1047480611"
$ws.Range("D16").Value = "DANNA MARCELA HERNANDEZ DIAZ"
$ws.Range("E16").Value = "2507"
$ws.Range("F16").Value = 1898
$ws.Range("G16").Value = 1423500

$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1007597829"
$ws.Range("D17").Value = "NOHELY JOHANA PEREIRA PAJARO"
$ws.Range("E17").Value = "2505"
$ws.Range("F17").Value = 1898
$ws.Range("G17").Value = 1423500

$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1071356875"
$ws.Range("D18").Value = "DERLIS LORED CALLE MUÑOZ"
$ws.Range("E18").Value = "2505"
$ws.Range("F18").Value = 1898
$ws.Range("G18").Value = 1423500

# Row 19 keeps the original worker (CC 1143381697 / ANA ELVIRA CHICO PADILLA / 2504 / 11388 / 1423500)
# already carried over by the row insert/shift above, so nothing else to set there.

# Adjust column D width to fit the longer names now present
$ws.Columns("D").ColumnWidth = 32.26953125
